$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet: insert a new (blank) column before column N,
# shifting the old Late / heading / Outstanding columns one place to the right.
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N").Insert()
$wsRepay.Columns("N").ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab and select cell J21 on it.
$wsRepay.Activate()
$wsRepay.Range("J21").Select()

Write-Output "done"
